$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 29 (SP10 / 9x21mm), shifting
# the rest of the ammo chart down by two rows.
$ws.Rows("29:30").Insert()

# Fill the newly inserted rows with the new ammo entries (9x19mm chart).
$ws.Range("A29").Value = "7n31"
$ws.Range("B29").Value = "9x19mm"
$ws.Range("C29").Value = 52
$ws.Range("D29").Value = 39

$ws.Range("A30").Value = "Quakemaker"
$ws.Range("B30").Value = "9x19mm"
$ws.Range("C30").Value = 85
$ws.Range("D30").Value = 8

# Match the saved view state (selected cell) captured in the workbook
# after the edit.
$ws.Range("D30").Select()
